$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.031606185696898
$ws.Cells.Item(2, 4).Value = 1.040962096974562
$ws.Cells.Item(2, 5).Value = 1.041161578228174
$ws.Cells.Item(2, 6).Value = 1.052119757537961
$ws.Cells.Item(2, 9).Value = 1.035219596958237
$ws.Cells.Item(2, 10).Value = 1.036740985899744
$ws.Cells.Item(2, 11).Value = 1.043742960062516
$ws.Cells.Item(2, 12).Value = 1.043941876243726
$ws.Cells.Item(2, 13).Value = 1.054869371120294
$ws.Cells.Item(2, 14).Value = 1.016236842976457

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.032493531424212
$ws.Cells.Item(3, 4).Value = 1.041662918549766
$ws.Cells.Item(3, 5).Value = 1.041946912461114
$ws.Cells.Item(3, 6).Value = 1.052983389296652
$ws.Cells.Item(3, 9).Value = 1.035376533910007
$ws.Cells.Item(3, 10).Value = 1.037270593465257
$ws.Cells.Item(3, 11).Value = 1.044254600510251
$ws.Cells.Item(3, 12).Value = 1.044537848008832
$ws.Cells.Item(3, 13).Value = 1.055545654497446
$ws.Cells.Item(3, 14).Value = 1.016412576592564

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.033068100219286
$ws.Cells.Item(4, 4).Value = 1.042116309068197
$ws.Cells.Item(4, 5).Value = 1.042455766023624
$ws.Cells.Item(4, 6).Value = 1.053542796666759
$ws.Cells.Item(4, 9).Value = 1.035476280028681
$ws.Cells.Item(4, 10).Value = 1.037613047267762
$ws.Cells.Item(4, 11).Value = 1.044584904017686
$ws.Cells.Item(4, 12).Value = 1.044923513636158
$ws.Cells.Item(4, 13).Value = 1.055983192204669
$ws.Cells.Item(4, 14).Value = 1.016526184417041

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.033309742567845
$ws.Cells.Item(5, 4).Value = 1.042306891648566
$ws.Cells.Item(5, 5).Value = 1.042669851531046
$ws.Cells.Item(5, 6).Value = 1.053778108552016
$ws.Cells.Item(5, 9).Value = 1.0355177808842
$ws.Cells.Item(5, 10).Value = 1.037756956744562
$ws.Cells.Item(5, 11).Value = 1.044723579778799
$ws.Cells.Item(5, 12).Value = 1.045085653802188
$ws.Cells.Item(5, 13).Value = 1.056167116457187
$ws.Cells.Item(5, 14).Value = 1.016573919922571

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.033350320832194
$ws.Cells.Item(6, 4).Value = 1.042338889914785
$ws.Cells.Item(6, 5).Value = 1.04270580696707
$ws.Cells.Item(6, 6).Value = 1.053817626435981
$ws.Cells.Item(6, 9).Value = 1.035524723678153
$ws.Cells.Item(6, 10).Value = 1.037781116337035
$ws.Cells.Item(6, 11).Value = 1.044746853228648
$ws.Cells.Item(6, 12).Value = 1.045112878182185
$ws.Cells.Item(6, 13).Value = 1.05619799712989
$ws.Cells.Item(6, 14).Value = 1.016581933432398

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.033071328689993
$ws.Cells.Item(7, 4).Value = 1.042118855733739
$ws.Cells.Item(7, 5).Value = 1.042458626004436
$ws.Cells.Item(7, 6).Value = 1.053545940379531
$ws.Cells.Item(7, 9).Value = 1.035476836265331
$ws.Cells.Item(7, 10).Value = 1.037614970423043
$ws.Cells.Item(7, 11).Value = 1.044586757734093
$ws.Cells.Item(7, 12).Value = 1.04492568013778
$ws.Cells.Item(7, 13).Value = 1.055985649876717
$ws.Cells.Item(7, 14).Value = 1.016526822360855

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.031905985367006
$ws.Cells.Item(8, 4).Value = 1.041198960058566
$ws.Cells.Item(8, 5).Value = 1.041426841778992
$ws.Cells.Item(8, 6).Value = 1.052411504789695
$ws.Cells.Item(8, 9).Value = 1.035273007250631
$ws.Cells.Item(8, 10).Value = 1.036920017757826
$ws.Cells.Item(8, 11).Value = 1.043916028305355
$ws.Cells.Item(8, 12).Value = 1.044143280255795
$ws.Cells.Item(8, 13).Value = 1.055097936338611
$ws.Cells.Item(8, 14).Value = 1.016296254164303

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029855600337107
$ws.Cells.Item(9, 4).Value = 1.039577390051027
$ws.Cells.Item(9, 5).Value = 1.039614062745056
$ws.Cells.Item(9, 6).Value = 1.050417004548839
$ws.Cells.Item(9, 9).Value = 1.034900067044837
$ws.Cells.Item(9, 10).Value = 1.035693654375954
$ws.Cells.Item(9, 11).Value = 1.042728341915022
$ws.Cells.Item(9, 12).Value = 1.042764895882822
$ws.Cells.Item(9, 13).Value = 1.053533258480599
$ws.Cells.Item(9, 14).Value = 1.015889189087978

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028490840226023
$ws.Cells.Item(10, 4).Value = 1.038496045458855
$ws.Cells.Item(10, 5).Value = 1.038409243453819
$ws.Cells.Item(10, 6).Value = 1.049090482714645
$ws.Cells.Item(10, 9).Value = 1.034642234712938
$ws.Cells.Item(10, 10).Value = 1.034874960298546
$ws.Cells.Item(10, 11).Value = 1.041932749597657
$ws.Cells.Item(10, 12).Value = 1.04184625611237
$ws.Cells.Item(10, 13).Value = 1.05248995339569
$ws.Cells.Item(10, 14).Value = 1.015617316837333

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.027900413606421
$ws.Cells.Item(11, 4).Value = 1.038027762150748
$ws.Cells.Item(11, 5).Value = 1.037888441629844
$ws.Cells.Item(11, 6).Value = 1.048516852165605
$ws.Cells.Item(11, 9).Value = 1.034528417583884
$ws.Cells.Item(11, 10).Value = 1.034520206040079
$ws.Cells.Item(11, 11).Value = 1.041587364696429
$ws.Cells.Item(11, 12).Value = 1.041448558506805
$ws.Cells.Item(11, 13).Value = 1.052038165568872
$ws.Cells.Item(11, 14).Value = 1.015499480862956

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.027681182812166
$ws.Cells.Item(12, 4).Value = 1.037853814481896
$ws.Cells.Item(12, 5).Value = 1.03769512843948
$ws.Cells.Item(12, 6).Value = 1.048303896506251
$ws.Cells.Item(12, 9).Value = 1.034485815022272
$ws.Cells.Item(12, 10).Value = 1.03438839745619
$ws.Cells.Item(12, 11).Value = 1.041458941264257
$ws.Cells.Item(12, 12).Value = 1.041300849150482
$ws.Cells.Item(12, 13).Value = 1.051870348479745
$ws.Cells.Item(12, 14).Value = 1.015455694698734

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.027728204923541
$ws.Cells.Item(13, 4).Value = 1.037891127109488
$ws.Cells.Item(13, 5).Value = 1.037736588603533
$ws.Cells.Item(13, 6).Value = 1.048349570928498
$ws.Cells.Item(13, 9).Value = 1.034494968147067
$ws.Cells.Item(13, 10).Value = 1.034416672503312
$ws.Cells.Item(13, 11).Value = 1.041486494482868
$ws.Cells.Item(13, 12).Value = 1.04133253268875
$ws.Cells.Item(13, 13).Value = 1.051906345925457
$ws.Cells.Item(13, 14).Value = 1.015465087725796

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027882290282352
$ws.Cells.Item(14, 4).Value = 1.038013383702831
$ws.Cells.Item(14, 5).Value = 1.037872459521502
$ws.Cells.Item(14, 6).Value = 1.048499246802379
$ws.Cells.Item(14, 9).Value = 1.034524902684645
$ws.Cells.Item(14, 10).Value = 1.034509311450716
$ws.Cells.Item(14, 11).Value = 1.041576751862583
$ws.Cells.Item(14, 12).Value = 1.041436348527692
$ws.Cells.Item(14, 13).Value = 1.052024293807902
$ws.Cells.Item(14, 14).Value = 1.015495861821952

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.027977237988887
$ws.Cells.Item(15, 4).Value = 1.038088709228197
$ws.Cells.Item(15, 5).Value = 1.037956192121919
$ws.Cells.Item(15, 6).Value = 1.048591482507382
$ws.Cells.Item(15, 9).Value = 1.034543303192123
$ws.Cells.Item(15, 10).Value = 1.034566384487686
$ws.Cells.Item(15, 11).Value = 1.041632344948251
$ws.Cells.Item(15, 12).Value = 1.041500314685954
$ws.Cells.Item(15, 13).Value = 1.05209696504196
$ws.Cells.Item(15, 14).Value = 1.015514820566128

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.028530036036113
$ws.Cells.Item(16, 4).Value = 1.038527122888891
$ws.Cells.Item(16, 5).Value = 1.038443826308616
$ws.Cells.Item(16, 6).Value = 1.049128568893283
$ws.Cells.Item(16, 9).Value = 1.034649742645993
$ws.Cells.Item(16, 10).Value = 1.034898498910491
$ws.Cells.Item(16, 11).Value = 1.041955653057429
$ws.Cells.Item(16, 12).Value = 1.041872651778328
$ws.Cells.Item(16, 13).Value = 1.052519936566028
$ws.Cells.Item(16, 14).Value = 1.015625134867078

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.028876932716106
$ws.Cells.Item(17, 4).Value = 1.038802114956963
$ws.Cells.Item(17, 5).Value = 1.038749946566869
$ws.Cells.Item(17, 6).Value = 1.049465674290393
$ws.Cells.Item(17, 9).Value = 1.034715927828151
$ws.Cells.Item(17, 10).Value = 1.035106758239078
$ws.Cells.Item(17, 11).Value = 1.042158218847941
$ws.Cells.Item(17, 12).Value = 1.042106231197576
$ws.Cells.Item(17, 13).Value = 1.052785248505926
$ws.Cells.Item(17, 14).Value = 1.015694302020794

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.029079322105676
$ws.Cells.Item(18, 4).Value = 1.038962507856792
$ws.Cells.Item(18, 5).Value = 1.038928587373882
$ws.Cells.Item(18, 6).Value = 1.049662375584582
$ws.Cells.Item(18, 9).Value = 1.034754322708828
$ws.Cells.Item(18, 10).Value = 1.035228207766693
$ws.Cells.Item(18, 11).Value = 1.042276286260787
$ws.Cells.Item(18, 12).Value = 1.042242481702317
$ws.Cells.Item(18, 13).Value = 1.05293999749169
$ws.Cells.Item(18, 14).Value = 1.015734635072906

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029148340171071
$ws.Cells.Item(19, 4).Value = 1.039017196736062
$ws.Cells.Item(19, 5).Value = 1.038989513816085
$ws.Cells.Item(19, 6).Value = 1.049729458042217
$ws.Cells.Item(19, 9).Value = 1.034767378757704
$ws.Cells.Item(19, 10).Value = 1.035269614706797
$ws.Cells.Item(19, 11).Value = 1.042316529605942
$ws.Cells.Item(19, 12).Value = 1.042288940849683
$ws.Cells.Item(19, 13).Value = 1.052992762345318
$ws.Cells.Item(19, 14).Value = 1.015748385726605

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.028839708760286
$ws.Cells.Item(20, 4).Value = 1.038772611448555
$ws.Cells.Item(20, 5).Value = 1.038717093848838
$ws.Cells.Item(20, 6).Value = 1.04942949846475
$ws.Cells.Item(20, 9).Value = 1.03470884848088
$ws.Cells.Item(20, 10).Value = 1.035084416521797
$ws.Cells.Item(20, 11).Value = 1.042136494319204
$ws.Cells.Item(20, 12).Value = 1.042081169542869
$ws.Cells.Item(20, 13).Value = 1.052756783342979
$ws.Cells.Item(20, 14).Value = 1.015686882172339

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027836913767546
$ws.Cells.Item(21, 4).Value = 1.037977382346177
$ws.Cells.Item(21, 5).Value = 1.037832445160794
$ws.Cells.Item(21, 6).Value = 1.048455167764421
$ws.Cells.Item(21, 9).Value = 1.034516096696946
$ws.Cells.Item(21, 10).Value = 1.034482032603387
$ws.Cells.Item(21, 11).Value = 1.041550176949804
$ws.Cells.Item(21, 12).Value = 1.041405776974693
$ws.Cells.Item(21, 13).Value = 1.051989561166311
$ws.Cells.Item(21, 14).Value = 1.015486800073587

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.027206879906832
$ws.Cells.Item(22, 4).Value = 1.03747735407185
$ws.Cells.Item(22, 5).Value = 1.03727701790515
$ws.Cells.Item(22, 6).Value = 1.047843240938209
$ws.Cells.Item(22, 9).Value = 1.034393021349787
$ws.Cells.Item(22, 10).Value = 1.034103076093559
$ws.Cells.Item(22, 11).Value = 1.041180772307618
$ws.Cells.Item(22, 12).Value = 1.04098120811561
$ws.Cells.Item(22, 13).Value = 1.051507161667335
$ws.Cells.Item(22, 14).Value = 1.015360904466469

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.027540828445229
$ws.Cells.Item(23, 4).Value = 1.037742431380906
$ws.Cells.Item(23, 5).Value = 1.037571385302652
$ws.Cells.Item(23, 6).Value = 1.048167570494539
$ws.Cells.Item(23, 9).Value = 1.034458444293364
$ws.Cells.Item(23, 10).Value = 1.034303987994553
$ws.Cells.Item(23, 11).Value = 1.041376672637164
$ws.Cells.Item(23, 12).Value = 1.041206272344546
$ws.Cells.Item(23, 13).Value = 1.05176289182727
$ws.Cells.Item(23, 14).Value = 1.015427653076148

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.028856528507941
$ws.Cells.Item(24, 4).Value = 1.038785942830665
$ws.Cells.Item(24, 5).Value = 1.038731938312323
$ws.Cells.Item(24, 6).Value = 1.049445844536523
$ws.Cells.Item(24, 9).Value = 1.034712047981555
$ws.Cells.Item(24, 10).Value = 1.03509451185832
$ws.Cells.Item(24, 11).Value = 1.042146310963815
$ws.Cells.Item(24, 12).Value = 1.042092493802575
$ws.Cells.Item(24, 13).Value = 1.05276964553421
$ws.Cells.Item(24, 14).Value = 1.015690234916436

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030385298566936
$ws.Cells.Item(25, 4).Value = 1.039996664790262
$ws.Cells.Item(25, 5).Value = 1.040082064615537
$ws.Cells.Item(25, 6).Value = 1.050932084072718
$ws.Cells.Item(25, 9).Value = 1.034998107275298
$ws.Cells.Item(25, 10).Value = 1.036010900973507
$ws.Cells.Item(25, 11).Value = 1.043036063294843
$ws.Cells.Item(25, 12).Value = 1.043121197117751
$ws.Cells.Item(25, 13).Value = 1.053937804749991
$ws.Cells.Item(25, 14).Value = 1.015994514345815
